$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("M2").Value = 1.532593
$ws.Range("N2").Value = 3.065186
$ws.Range("O2").Value = 0.0795983245703594
$ws.Range("P2").Value = 0.05838920196386116
$ws.Range("Q2").Value = 0.1859318838705
$ws.Range("R2").Value = 0.743727535482
$ws.Range("S2").Value = 0.0795983245703594
$ws.Range("T2").Value = 0.05838920196386116

# Row 3
$ws.Range("O3").Value = 0.6326044366842063
$ws.Range("P3").Value = 0.6960687002426557
$ws.Range("S3").Value = 0.6326044366842063
$ws.Range("T3").Value = 0.6960687002426557

# Row 4
$ws.Range("M4").Value = 0.8528209999999999
$ws.Range("N4").Value = 2.558463
$ws.Range("O4").Value = 0.04429298760885536
$ws.Range("P4").Value = 0.04873655720209673
$ws.Range("Q4").Value = 0.1034629644885
$ws.Range("R4").Value = 0.6207777869309999
$ws.Range("S4").Value = 0.04429298760885536
$ws.Range("T4").Value = 0.04873655720209673

# Row 5
$ws.Range("M5").Value = 3.7338975
$ws.Range("N5").Value = 7.467795
$ws.Range("O5").Value = 0.1939275366111247
$ws.Range("P5").Value = 0.142255181408147
$ws.Range("Q5").Value = 0.45299084385375
$ws.Range("R5").Value = 1.811963375415
$ws.Range("S5").Value = 0.1939275366111247
$ws.Range("T5").Value = 0.142255181408147

# Row 6
$ws.Range("M6").Value = 0.2147316666666667
$ws.Range("N6").Value = 0.644195
$ws.Range("O6").Value = 0.01115252444639089
$ws.Range("P6").Value = 0.01227137014168456
$ws.Range("Q6").Value = 0.0260509237025
$ws.Range("R6").Value = 0.156305542215
$ws.Range("S6").Value = 0.01115252444639089
$ws.Range("T6").Value = 0.01227137014168456

# Row 7
$ws.Range("M7").Value = 0.7398226666666666
$ws.Range("N7").Value = 2.219468
$ws.Range("O7").Value = 0.03842419007906348
$ws.Range("P7").Value = 0.04227898904155473
$ws.Range("Q7").Value = 0.089754176186
$ws.Range("R7").Value = 0.538525057116
$ws.Range("S7").Value = 0.03842419007906348
$ws.Range("T7").Value = 0.04227898904155473
